function Set-NumValue($rng, $val) {
    # Columns L and M are formatted as Text (numFmtId 49 / "@"). Assigning a
    # numeric .Value straight to such a cell stores it as a text shared
    # string instead of a real number, so temporarily switch to a generic
    # number format, write the value, then restore the original format.
    $fmt = $rng.NumberFormat
    $rng.NumberFormat = "General"
    $rng.Value = $val
    $rng.NumberFormat = $fmt
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Title text update (A1): new "data as of" date ---
$ws.Range("A1").Value = "Données COVID-19 Valais 10.06.2020"

# --- 2. New day's figures arrived for 2020-06-09 (row 105): one new case ---
$ws.Range("C105").Value = 1

# --- 3. A death recorded retroactively on 2020-05-11 (row 76) ---
Set-NumValue $ws.Range("M76") 1

# --- 4. Insert the new "current" bottom row (106), carrying forward the
#        running totals with no new data yet for the very latest day ---

# 4a. Put the formulas in place on row 106 FIRST, while it is still blank
#     (default/General format) - this avoids the engine inferring a Text
#     number format for the SUM-like formulas from their (Text-formatted)
#     operands once the real borders/format get pasted in afterwards.
$ws.Range("B106").Formula = "=B105+C106"
$ws.Range("H106").Formula = "=G106+E106"
$ws.Range("J106").Formula = "=J105+K106"
$ws.Range("K106").Formula = "=L106+M106"

# 4b. Duplicate row 105's current ("final row" / thicker bottom border)
#     formatting down onto the new row 106.
$ws.Range("A105:M105").Copy()
$ws.Range("A106:M106").PasteSpecial(-4122)

# 4c. Demote row 105 back to a normal (non-final) row by copying the
#     normal formatting from row 104 onto it.
$ws.Range("A104:M104").Copy()
$ws.Range("A105:M105").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# 4d. Populate row 106's plain values (date left blank - no new day yet).
$ws.Range("A106").Value = ""
$ws.Range("C106").Value = 0
$ws.Range("D106").Value = 0
$ws.Range("E106").Value = 4
$ws.Range("F106").Value = 3
$ws.Range("G106").Value = 11
$ws.Range("I106").Value = 0
Set-NumValue $ws.Range("L106") 0
Set-NumValue $ws.Range("M106") 0

# --- 5. Update the view: scroll near the bottom and select the new row ---
$win = $excel.ActiveWindow
$win.ScrollRow = 88
$win.ScrollColumn = 1
$ws.Range("A106").Select()
